$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 223, shifting existing rows 223:303 down to 224:304
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new data record
$ws.Range("A223").Value = 6
$ws.Range("B223").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C223").Value = "Metropolitana"
$ws.Range("D223").Value = 44553
$ws.Range("E223").Value = 13
$ws.Range("F223").Value = 100112032
$ws.Range("G223").Value = "Zapallo italiano"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 350
$ws.Range("K223").Value = 6000
$ws.Range("L223").Value = 7000
$ws.Range("M223").Value = 6343
$ws.Range("N223").Value = "`$/caja 50 unidades"
$ws.Range("O223").Value = "Región de O'Higgins"
$ws.Range("P223").Value = 127
$ws.Range("Q223").Value = 50
$ws.Range("R223").Value = "Hortaliza"
